$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill AM5:AN54 (test MSE columns for the "beef" prodnorm run) with the newly
# computed values so the AM56/AN56 (AVERAGE) and AM57/AN57 (STDEV.S) formulas
# stop showing #DIV/0! and recalculate to the expected statistics.
$ws.Cells.Item(5, 39).Value = 0.00006794603094038043
$ws.Cells.Item(5, 40).Value = 0.00034060974658983918
$ws.Cells.Item(6, 39).Value = 0.00006075710643024922
$ws.Cells.Item(6, 40).Value = 0.00081203189827220397
$ws.Cells.Item(7, 39).Value = 0.00005963859013682863
$ws.Cells.Item(7, 40).Value = 0.00097125230494909259
$ws.Cells.Item(8, 39).Value = 0.00008243796025329292
$ws.Cells.Item(8, 40).Value = 0.00037117552260130257
$ws.Cells.Item(9, 39).Value = 0.00008418906003950633
$ws.Cells.Item(9, 40).Value = 0.0005061530916661255
$ws.Cells.Item(10, 39).Value = 0.00007825910569935644
$ws.Cells.Item(10, 40).Value = 0.00031922570655601071
$ws.Cells.Item(11, 39).Value = 0.00007686877353164497
$ws.Cells.Item(11, 40).Value = 0.00043958391643057738
$ws.Cells.Item(12, 39).Value = 0.0000644623232538483
$ws.Cells.Item(12, 40).Value = 0.00078444376074833348
$ws.Cells.Item(13, 39).Value = 0.00006957469285668557
$ws.Cells.Item(13, 40).Value = 0.00067763822014559809
$ws.Cells.Item(14, 39).Value = 0.00006591351900090224
$ws.Cells.Item(14, 40).Value = 0.00086206526875943286
$ws.Cells.Item(15, 39).Value = 0.00007848073013479225
$ws.Cells.Item(15, 40).Value = 0.00032429946155577489
$ws.Cells.Item(16, 39).Value = 0.00009040992510164936
$ws.Cells.Item(16, 40).Value = 0.00038209616074169947
$ws.Cells.Item(17, 39).Value = 0.00007135034482443704
$ws.Cells.Item(17, 40).Value = 0.00057984856852486976
$ws.Cells.Item(18, 39).Value = 0.00007762724145826452
$ws.Cells.Item(18, 40).Value = 0.00045226008976046068
$ws.Cells.Item(19, 39).Value = 0.00005703398751351211
$ws.Cells.Item(19, 40).Value = 0.00109039015626093702
$ws.Cells.Item(20, 39).Value = 0.00007155834809736865
$ws.Cells.Item(20, 40).Value = 0.00033824800008620251
$ws.Cells.Item(21, 39).Value = 0.00008131232674440621
$ws.Cells.Item(21, 40).Value = 0.00045522688709077332
$ws.Cells.Item(22, 39).Value = 0.0000776068708733972
$ws.Cells.Item(22, 40).Value = 0.00034326949868271448
$ws.Cells.Item(23, 39).Value = 0.00006852287877141553
$ws.Cells.Item(23, 40).Value = 0.00056308429688425507
$ws.Cells.Item(24, 39).Value = 0.00008421468399557183
$ws.Cells.Item(24, 40).Value = 0.00050682583087871557
$ws.Cells.Item(25, 39).Value = 0.00007861447261731927
$ws.Cells.Item(25, 40).Value = 0.0002043993733497789
$ws.Cells.Item(26, 39).Value = 0.00008210921176833968
$ws.Cells.Item(26, 40).Value = 0.00044502441381264932
$ws.Cells.Item(27, 39).Value = 0.00007940052463144631
$ws.Cells.Item(27, 40).Value = 0.00031238294625978103
$ws.Cells.Item(28, 39).Value = 0.00006023149500239734
$ws.Cells.Item(28, 40).Value = 0.00099713461616780016
$ws.Cells.Item(29, 39).Value = 0.00007704974914771672
$ws.Cells.Item(29, 40).Value = 0.00034553660402215111
$ws.Cells.Item(30, 39).Value = 0.00006974650548725389
$ws.Cells.Item(30, 40).Value = 0.00083533502519707195
$ws.Cells.Item(31, 39).Value = 0.00008339557974476773
$ws.Cells.Item(31, 40).Value = 0.00042961543597021099
$ws.Cells.Item(32, 39).Value = 0.00007977194335191365
$ws.Cells.Item(32, 40).Value = 0.00045404133677457462
$ws.Cells.Item(33, 39).Value = 0.00008130889914749843
$ws.Cells.Item(33, 40).Value = 0.00026800132200669992
$ws.Cells.Item(34, 39).Value = 0.00007083507558735815
$ws.Cells.Item(34, 40).Value = 0.00051190832175078947
$ws.Cells.Item(35, 39).Value = 0.0000761577185902921
$ws.Cells.Item(35, 40).Value = 0.00025684863983024191
$ws.Cells.Item(36, 39).Value = 0.00007378454681552652
$ws.Cells.Item(36, 40).Value = 0.00083199655750667006
$ws.Cells.Item(37, 39).Value = 0.00007951335173658155
$ws.Cells.Item(37, 40).Value = 0.00036820137389160091
$ws.Cells.Item(38, 39).Value = 0.00005924833516053748
$ws.Cells.Item(38, 40).Value = 0.00068394518233425711
$ws.Cells.Item(39, 39).Value = 0.00007689140397135454
$ws.Cells.Item(39, 40).Value = 0.00072523234244959507
$ws.Cells.Item(40, 39).Value = 0.00008003707863546856
$ws.Cells.Item(40, 40).Value = 0.00034715119422411799
$ws.Cells.Item(41, 39).Value = 0.00007445051686098719
$ws.Cells.Item(41, 40).Value = 0.00033358256708779813
$ws.Cells.Item(42, 39).Value = 0.00007809521429764983
$ws.Cells.Item(42, 40).Value = 0.00043077227979823482
$ws.Cells.Item(43, 39).Value = 0.00007382181653606647
$ws.Cells.Item(43, 40).Value = 0.0003881472386070956
$ws.Cells.Item(44, 39).Value = 0.00006478645412950882
$ws.Cells.Item(44, 40).Value = 0.00070328697446350545
$ws.Cells.Item(45, 39).Value = 0.00007416251528002859
$ws.Cells.Item(45, 40).Value = 0.0003850192093970017
$ws.Cells.Item(46, 39).Value = 0.00007792990729445719
$ws.Cells.Item(46, 40).Value = 0.00039305080870227169
$ws.Cells.Item(47, 39).Value = 0.00008325057650912358
$ws.Cells.Item(47, 40).Value = 0.00031128766317533018
$ws.Cells.Item(48, 39).Value = 0.00007846028678312144
$ws.Cells.Item(48, 40).Value = 0.00062847719840456091
$ws.Cells.Item(49, 39).Value = 0.00008404107600267388
$ws.Cells.Item(49, 40).Value = 0.00058310318907220577
$ws.Cells.Item(50, 39).Value = 0.00006789461475644394
$ws.Cells.Item(50, 40).Value = 0.00094862891141102561
$ws.Cells.Item(51, 39).Value = 0.00007687729432133395
$ws.Cells.Item(51, 40).Value = 0.00050732335540102096
$ws.Cells.Item(52, 39).Value = 0.00006524528539157814
$ws.Cells.Item(52, 40).Value = 0.0004963573287806686
$ws.Cells.Item(53, 39).Value = 0.00007617857538505319
$ws.Cells.Item(53, 40).Value = 0.00036274795197884641
$ws.Cells.Item(54, 39).Value = 0.00007618951326089121
$ws.Cells.Item(54, 40).Value = 0.00063140706335729462

# Recalculate so AM56/AN56/AM57/AN57 formulas update from #DIV/0! to real values
$excel.CalculateFullRebuild()

# Update the view state (scrolled position + active selection) to match the saved workbook
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 7
$win.Zoom = 70
$ws.Range("AM61").Select()
